$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-blank column B values across the cash flow rows
$ws.Range("B3").Value = 370602000.0
$ws.Range("B4").Value = -182600000.0
$ws.Range("B5").Value = 22169000.0
$ws.Range("B6").Value = 602904000.0
$ws.Range("B7").Value = -479244000.0
$ws.Range("F7").Value = -378075000.0
$ws.Range("B8").Value = -535289000.0
$ws.Range("C8").Value = -534287000.0
$ws.Range("D8").Value = -382712000.0
$ws.Range("F8").Value = -405448000.0
$ws.Range("B9").Value = 94408000.0
$ws.Range("B10").Value = -908696000.0
$ws.Range("B11").Value = 305743000.0
$ws.Range("B12").Value = 212772000.0
$ws.Range("B14").Value = -365127000.0
$ws.Range("B15").Value = -15868000.0
$ws.Range("B16").Value = 137520000.0
$ws.Range("B17").Value = -168272000.0
$ws.Range("B18").Value = 67223000.0
$ws.Range("B19").Value = -101049000.0
$ws.Range("B20").Value = -367691000.0
$ws.Range("B21").Value = 50722000.0
$ws.Range("B23").Value = 212772000.0
$ws.Range("B24").Value = 212772000.0
